# Commit: "case with 380 kV done"
# Updates the vm_pu (per-unit voltage magnitude) results table for rows 2-25.
# Columns B:F and I:N are updated per row; columns G and H are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$endRow = 25

# New values for columns B,C,D,E,F for each data row (row 2 first)
$bfData = @(
    @(1.02,1.048424262002884,1.050466937145047,1.061821906073496,1.068790066389022),
    @(1.02,1.049462215234939,1.051275161843366,1.062882318556933,1.069991803287285),
    @(1.02,1.050133906502342,1.051798071366729,1.063568942066999,1.070770173041282),
    @(1.02,1.050416301513871,1.052017886288723,1.063857710015529,1.071097583147129),
    @(1.02,1.050463717822441,1.052054793215427,1.063906201956747,1.071152567507862),
    @(1.02,1.050137679815202,1.051801008608879,1.063572800160721,1.070774547190042),
    @(1.02,1.048775029570439,1.050740093832206,1.062180180962782,1.069196040574343),
    @(1.02,1.046374363984585,1.048870141849277,1.059729777569696,1.066420365995713),
    @(1.02,1.044774240967651,1.047623206370883,1.058098575892874,1.064573839087166),
    @(1.02,1.044081442011265,1.047083202765532,1.057392816102861,1.063775198840783),
    @(1.02,1.043824114916129,1.046882610802623,1.057130749639863,1.063478685760139),
    @(1.02,1.043879312038813,1.04692563893382,1.057186959980863,1.063542282609295),
    @(1.02,1.044060171086574,1.047066621985497,1.057371151891752,1.063750686175454),
    @(1.02,1.044171605531481,1.047153484903164,1.057484649752025,1.063879108692476),
    @(1.02,1.044820221079951,1.047659043132415,1.058145426656151,1.064626861513294),
    @(1.02,1.045227097518099,1.047976147293486,1.058560064697899,1.065096152571446),
    @(1.02,1.045464427950443,1.048161101789368,1.058801970222028,1.065369970719436),
    @(1.02,1.045545352575639,1.048224165297447,1.058884462947997,1.065463350676422),
    @(1.02,1.045183442879539,1.04794212573956,1.058515572358297,1.065045792919736),
    @(1.02,1.044006912352208,1.047025106280476,1.057316909643443,1.063689312728529),
    @(1.02,1.043267234402012,1.046448478676134,1.056563749579306,1.062837234902784),
    @(1.02,1.043659346761986,1.046754165573914,1.056962967900841,1.06328886226797),
    @(1.02,1.045203168511953,1.047957498636884,1.058535676363866,1.065068548004601),
    @(1.02,1.046994936016001,1.049353624104693,1.060362842510848,1.067137251947214)
)

# New values for columns I,J,K,L,M,N for each data row (row 2 first)
$inData = @(
    @(1.046123027601769,1.053468087355162,1.05322113440855,1.064544924738974,1.071494306675087,1.054964132873685),
    @(1.04644136644952,1.054154478011151,1.053841822157775,1.065419414454857,1.072511131815288,1.055651498283135),
    @(1.046646015278282,1.054598085243874,1.054242730077543,1.065985130371908,1.073169274899656,1.056095735489012),
    @(1.046731729281459,1.05478444990666,1.054411099521599,1.066222924399436,1.073446003622453,1.056282364810995),
    @(1.046746102250023,1.054815733854946,1.054439359385155,1.066262849129781,1.073492470255293,1.056313693186075),
    @(1.046647161851947,1.054600575960002,1.05424498051538,1.065988307916848,1.073172972384029,1.056098229742241),
    @(1.046230888612707,1.053700166669133,1.053431047388879,1.064840491053614,1.071837908352729,1.055196541766902),
    @(1.045487121857094,1.052109449866972,1.051991299342667,1.062816840479855,1.069486800678214,1.053603565964424),
    @(1.044984397864271,1.051046232654669,1.051027784218048,1.061467033411151,1.067920359777742,1.052538838861727),
    @(1.044765081779776,1.050585198287238,1.050609699540988,1.060882383487909,1.067242299621525,1.052077149772613),
    @(1.04468337262311,1.050413851212832,1.050454272667067,1.060665192064854,1.066990470469393,1.051905559365708),
    @(1.0447009106196,1.050450610195277,1.050487618214101,1.060711781560438,1.067044487208404,1.051942370550123),
    @(1.044758332681526,1.050571036692065,1.050596854591668,1.060864430916947,1.067221482682548,1.052062968066349),
    @(1.044793679790992,1.050645222392209,1.050664141290792,1.060958479718466,1.067330539743005,1.052137259118697),
    @(1.044998918737785,1.051076816194931,1.051055512662544,1.061505831007546,1.067965364984479,1.052569465834124),
    @(1.045127222420043,1.051347368352672,1.051300775010142,1.061849123146057,1.068363632506152,1.052840402206961),
    @(1.04520190220054,1.051505113724046,1.051443747837962,1.06204934281271,1.06859595637484,1.052998371594857),
    @(1.045227339346879,1.051558890114611,1.051492483535133,1.062117609695941,1.068675176409858,1.053052224354065),
    @(1.045113472949174,1.05131834719642,1.051274469426413,1.061812292890474,1.068320899983703,1.052811339837339),
    @(1.044741430089801,1.05053557680681,1.050564690833866,1.060819480216664,1.067169361030208,1.052027457823991),
    @(1.044506092174586,1.05004284874261,1.050117663428998,1.060195106593392,1.066445531122301,1.051534030029423),
    @(1.044630983842215,1.050304107197813,1.050354713275984,1.060526113418329,1.066829229311096,1.051795659501598),
    @(1.045119686230058,1.051331460802646,1.051286356048533,1.06182893494947,1.06834020890561,1.05282447206639),
    @(1.045680616346113,1.052521171427327,1.052364158674508,1.063340127197997,1.070094447852153,1.05401587221662)
)

$bfCols = @(2,3,4,5,6)      # B,C,D,E,F
$inCols = @(9,10,11,12,13,14) # I,J,K,L,M,N

for ($r = $startRow; $r -le $endRow; $r++) {
    $rowIdx = $r - $startRow
    $bfRow = $bfData[$rowIdx]
    for ($i = 0; $i -lt $bfCols.Length; $i++) {
        $ws.Cells.Item($r, $bfCols[$i]).Value = $bfRow[$i]
    }
    $inRow = $inData[$rowIdx]
    for ($i = 0; $i -lt $inCols.Length; $i++) {
        $ws.Cells.Item($r, $inCols[$i]).Value = $inRow[$i]
    }
}
